$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1702127659574468
$ws.Range("C2").Value = 0.6063829787234043
$ws.Range("J2").Value = 0.01063829787234043
$ws.Range("P2").Value = 0.1347517730496454
$ws.Range("S2").Value = 0.07801418439716312
$ws.Range("C3").Value = 0.03314917127071823
$ws.Range("J3").Value = 0.02762430939226519
$ws.Range("P3").Value = 0.7624309392265194
$ws.Range("S3").Value = 0.1767955801104972
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.7708333333333334
$ws.Range("S4").Value = 0.1875
$ws.Range("B6").Value = 0.05825242718446602
$ws.Range("D6").Value = 0.01456310679611651
$ws.Range("F6").Value = 0.0825242718446602
$ws.Range("J6").Value = 0.2427184466019418
$ws.Range("O6").Value = 0.004854368932038835
$ws.Range("Q6").Value = 0.1213592233009709
$ws.Range("R6").Value = 0.09223300970873786
$ws.Range("S6").Value = 0.383495145631068
$ws.Range("B7").Value = 0.08125
$ws.Range("D7").Value = 0.03125
$ws.Range("F7").Value = 0.06875000000000001
$ws.Range("J7").Value = 0.08125
$ws.Range("O7").Value = 0.00625
$ws.Range("Q7").Value = 0.18125
$ws.Range("R7").Value = 0.075
$ws.Range("S7").Value = 0.475
$ws.Range("B8").Value = 0.09302325581395349
$ws.Range("D8").Value = 0.02093023255813953
$ws.Range("F8").Value = 0.08604651162790698
$ws.Range("J8").Value = 0.1116279069767442
$ws.Range("O8").Value = 0.02325581395348837
$ws.Range("Q8").Value = 0.1488372093023256
$ws.Range("R8").Value = 0.1046511627906977
$ws.Range("S8").Value = 0.4116279069767442
$ws.Range("B9").Value = 0.1420118343195266
$ws.Range("D9").Value = 0.01775147928994083
$ws.Range("F9").Value = 0.08284023668639054
$ws.Range("J9").Value = 0.1005917159763314
$ws.Range("O9").Value = 0.01183431952662722
$ws.Range("Q9").Value = 0.1834319526627219
$ws.Range("R9").Value = 0.08284023668639054
$ws.Range("S9").Value = 0.378698224852071
$ws.Range("B10").Value = 0.1293800539083558
$ws.Range("D10").Value = 0.02515723270440252
$ws.Range("E10").Value = 0.0008984725965858042
$ws.Range("F10").Value = 0.0664869721473495
$ws.Range("J10").Value = 0.1114106019766397
$ws.Range("O10").Value = 0.01347708894878706
$ws.Range("Q10").Value = 0.1949685534591195
$ws.Range("R10").Value = 0.09344115004492363
$ws.Range("S10").Value = 0.3647798742138365
$ws.Range("G11").Value = 0.1451612903225807
$ws.Range("J11").Value = 0.1088709677419355
$ws.Range("K11").Value = 0.2096774193548387
$ws.Range("L11").Value = 0.5282258064516129
$ws.Range("S11").Value = 0.008064516129032258
$ws.Range("G12").Value = 0.7686567164179104
$ws.Range("J12").Value = 0.1567164179104478
$ws.Range("K12").Value = 0.007462686567164179
$ws.Range("L12").Value = 0.02985074626865672
$ws.Range("S12").Value = 0.03731343283582089
$ws.Range("G13").Value = 0.675
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.025
$ws.Range("F15").Value = 0.01142857142857143
$ws.Range("H15").Value = 0.2514285714285714
$ws.Range("I15").Value = 0.06857142857142857
$ws.Range("J15").Value = 0.32
$ws.Range("K15").Value = 0.04
$ws.Range("M15").Value = 0.01142857142857143
$ws.Range("O15").Value = 0.04571428571428571
$ws.Range("S15").Value = 0.2514285714285714
$ws.Range("F16").Value = 0.02898550724637681
$ws.Range("H16").Value = 0.178743961352657
$ws.Range("I16").Value = 0.07246376811594203
$ws.Range("J16").Value = 0.4057971014492754
$ws.Range("K16").Value = 0.106280193236715
$ws.Range("M16").Value = 0.03381642512077294
$ws.Range("O16").Value = 0.04347826086956522
$ws.Range("S16").Value = 0.1304347826086956
$ws.Range("F17").Value = 0.01639344262295082
$ws.Range("H17").Value = 0.1639344262295082
$ws.Range("I17").Value = 0.09836065573770492
$ws.Range("J17").Value = 0.4699453551912569
$ws.Range("K17").Value = 0.07650273224043716
$ws.Range("M17").Value = 0.01366120218579235
$ws.Range("N17").Value = 0.00273224043715847
$ws.Range("O17").Value = 0.06830601092896176
$ws.Range("S17").Value = 0.09016393442622951
$ws.Range("F18").Value = 0.02094240837696335
$ws.Range("H18").Value = 0.1884816753926702
$ws.Range("I18").Value = 0.07853403141361257
$ws.Range("J18").Value = 0.418848167539267
$ws.Range("K18").Value = 0.09947643979057591
$ws.Range("M18").Value = 0.01047120418848168
$ws.Range("O18").Value = 0.07329842931937172
$ws.Range("S18").Value = 0.1099476439790576
$ws.Range("F19").Value = 0.01263537906137184
$ws.Range("H19").Value = 0.2310469314079422
$ws.Range("I19").Value = 0.08212996389891697
$ws.Range("J19").Value = 0.3763537906137184
$ws.Range("K19").Value = 0.1028880866425993
$ws.Range("M19").Value = 0.02256317689530686
$ws.Range("N19").Value = 0.0009025270758122744
$ws.Range("O19").Value = 0.05776173285198556
$ws.Range("S19").Value = 0.1137184115523466
